$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($r = 2; $r -le 121; $r++) {
    $c = $ws.Cells.Item($r, 1)
    $d = $c.Value()
    $yyyymmdd = [int]$d.ToString("yyyyMMdd")
    $c.Style = "Normal"
    $c.Value = $yyyymmdd
}
